$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-28
$values = @{
    2  = @(6, 6)
    3  = @(8, 8)
    4  = @(11, 11)
    5  = @(8, 9)
    6  = @(3, 5)
    7  = @(8, 9)
    8  = @(8, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(5, 6)
    13 = @(8, 9)
    14 = @(9, 9)
    15 = @(4, 6)
    16 = @(10, 10)
    17 = @(8, 8)
    18 = @(9, 9)
    19 = @(8, 9)
    20 = @(9, 10)
    21 = @(6, 6)
    22 = @(8, 8)
    23 = @(7, 7)
    24 = @(6, 6)
    25 = @(7, 7)
    26 = @(9, 9)
    27 = @(7, 7)
    28 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
